$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6045094728469849
$ws.Range("B1").Value = 1.482218146324158
$ws.Range("C1").Value = 5.823262691497803
$ws.Range("D1").Value = 2.178108930587769
$ws.Range("E1").Value = 1.423992156982422
